# Om Opplandskongene.docx edit
#
# The commit merges several "-" / text (and other split) w:r runs that
# were previously separate runs within the same paragraph into a single
# run per paragraph (no visible text change - only run consolidation),
# fixes up the section grid settings, and tweaks the "Normal" style's
# default paragraph properties in styles.xml.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge runs that are split across multiple <w:r> elements within the
#    same paragraph back into a single run, by re-finding & replacing
#    the paragraph's own (already-correct) text scoped to that
#    paragraph's Range. Word's Find/Replace naturally collapses the
#    matched span into a single run when the replacement text is
#    supplied literally (wdReplaceAll / match-type = 2).
# ---------------------------------------------------------------------

$merges = @(
    @{ Index = 7;  Text = "-Ryddet Värmland" },
    @{ Index = 8;  Text = "-Flyktet for Ivar Vidfavne" },
    @{ Index = 9;  Text = "-Var konge I Värmland til han ble “eldgammal”." },
    @{ Index = 10; Text = "-Kona het Sølva, søster til Sølve den gamle som rydda Solør (hedmark)" },
    @{ Index = 11; Text = "-Fikk to sønner sammen, Ingjald og Halfdan" },
    @{ Index = 17; Text = "-Halfdan giftet seg med Åsa, dattrer av kong Eysteinn Illråde av Heid" },
    @{ Index = 19; Text = "-Halfdan hersket over Solør, Raumariki og myte av heidmørk" },
    @{ Index = 25; Text = "-Agnar Eiriks far var sølnn av kong Sigtrygg av Vindli" },
    @{ Index = 26; Text = "-Halfdan var sønnen hans" },
    @{ Index = 27; Text = "-Døde av drukning" },
    @{ Index = 38; Text = "-Cynical/Trusting? Cruel" },
    @{ Index = 40; Text = "-Første kona var Alfarin av Alfheim og arva halve Vingulmark." },
    @{ Index = 41; Text = "-Sønnen deres var Olaf" },
    @{ Index = 43; Text = "-De hadde to sønner, Halfdan og Olaf" },
    @{ Index = 44; Text = "-Drept på ordre av Åsa sopm hevn for drapet på faren og sønnen Gyrd" },
    @{ Index = 47; Text = "Gjeveste, sterkeste og gildeste. Derfor kalt Geirstadalf" }
)

foreach ($m in $merges) {
    $p = $d.Paragraphs.Item($m.Index)
    $r = $p.Range
    $r.Find.Execute($m.Text, $false, $false, $false, $false, $false, $true, 1, $false, $m.Text, 2) | Out-Null
}

Write-Output "done"
